# Update countries & provincias Spain
# Refreshes the COVID country table ("Pais" sheet): bumps the "last updated"
# timestamp, updates the per-country counters, and re-sorts a handful of
# country rows whose totals crossed each other (the name + figures for the
# two countries involved simply trade places between the two adjacent rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 12:05"

# India (row 12)
$ws.Range("B12").Value = 182889
$ws.Range("C12").Value = 1062
$ws.Range("E12").Value = 90654

# Belgica (row 22)
$ws.Range("B22").Value = 58381
$ws.Range("C22").Value = 195
$ws.Range("D22").Value = 15887
$ws.Range("E22").Value = 33027
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 9467

# Indonesia overtakes Kuwait -> rows 35/36 swap name + figures
$ws.Range("A35").Value = "Indonesia"
$ws.Range("B35").Value = 26473
$ws.Range("C35").Value = 700
$ws.Range("D35").Value = 7308
$ws.Range("E35").Value = 17552
$ws.Range("G35").Value = 40
$ws.Range("H35").Value = 1613

$ws.Range("A36").Value = "Kuwait"
$ws.Range("B36").Value = 26192
$ws.Range("D36").Value = 10156
$ws.Range("E36").Value = 15831
$ws.Range("H36").Value = 205

# Rumania (row 41)
$ws.Range("B41").Value = 19257
$ws.Range("C41").Value = 124
$ws.Range("D41").Value = 13256
$ws.Range("E41").Value = 4739

# Barein overtakes Kazajistan -> rows 53/54 swap name + figures
$ws.Range("A53").Value = "Barein"
$ws.Range("B53").Value = 11288
$ws.Range("C53").Value = 495
$ws.Range("D53").Value = 6673
$ws.Range("E53").Value = 4597
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 18

$ws.Range("A54").Value = "Kazajistan"
$ws.Range("B54").Value = 10858
$ws.Range("C54").Value = 476
$ws.Range("D54").Value = 5220
$ws.Range("E54").Value = 5600
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 38

# Malasia overtakes Marruecos -> rows 64/65 swap name + figures
$ws.Range("A64").Value = "Malasia"
$ws.Range("B64").Value = 7819
$ws.Range("C64").Value = 57
$ws.Range("D64").Value = 6353
$ws.Range("E64").Value = 1351
$ws.Range("H64").Value = 115

$ws.Range("A65").Value = "Marruecos"
$ws.Range("B65").Value = 7783
$ws.Range("C65").Value = 3
$ws.Range("D65").Value = 5412
$ws.Range("E65").Value = 2167
$ws.Range("H65").Value = 204

# Finlandia (row 67)
$ws.Range("B67").Value = 6859
$ws.Range("C67").Value = 33
$ws.Range("E67").Value = 1039
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 320

# Consejo Danes para los Refugiados (row 82)
$ws.Range("E82").Value = 2550
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 72

# Eslovenia (row 104)
$ws.Range("D104").Value = 1358
$ws.Range("E104").Value = 7

# Albania (row 111)
$ws.Range("B111").Value = 1136
$ws.Range("C111").Value = 14
$ws.Range("D111").Value = 872
$ws.Range("E111").Value = 231

# Tunez (row 113)
$ws.Range("B113").Value = 1077
$ws.Range("C113").Value = 1
$ws.Range("D113").Value = 960
$ws.Range("E113").Value = 69

# Belice overtakes Santa Lucia -> rows 200/201 swap name + figures
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# Islas Virgenes Britanicas overtakes Papua Nueva Guinea -> rows 213/214 swap name + figures
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
